# Updates crypto price/volume data to the latest scraped values.
# Rows 17/18 and 31/32 had their underlying coin data swapped (re-ranked),
# so Coin (B), Link (C), Price (D) and Volume (E) are all replaced there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "61.545.38"
$ws.Range("E2").Value = "  +0.75%  "

# Row 3
$ws.Range("D3").Value = "2.938.35"
$ws.Range("E3").Value = "  +0.30%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.05%  "

# Row 5
$ws.Range("D5").Value = "'597.89"
$ws.Range("E5").Value = "  +1.11%  "

# Row 6
$ws.Range("D6").Value = "'145.34"
$ws.Range("E6").Value = "  -0.56%  "

# Row 7
$ws.Range("E7").Value = "  +0.07%  "

# Row 8
$ws.Range("D8").Value = "'0.503"
$ws.Range("E8").Value = "  -0.64%  "

# Row 9
$ws.Range("D9").Value = "'7.00"
$ws.Range("E9").Value = "  +2.04%  "

# Row 10
$ws.Range("E10").Value = "  -1.41%  "

# Row 11
$ws.Range("E11").Value = "  -0.28%  "

# Row 13
$ws.Range("D13").Value = "'33.64"
$ws.Range("E13").Value = "  -0.38%  "

# Row 14
$ws.Range("E14").Value = "  +0.57%  "

# Row 15
$ws.Range("D15").Value = "3.426.41"
$ws.Range("E15").Value = "  +0.44%  "

# Row 16
$ws.Range("D16").Value = "61.544.19"
$ws.Range("E16").Value = "  +0.83%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.936.51"
$ws.Range("E17").Value = "  +0.28%  "

# Row 18
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'6.70"
$ws.Range("E18").Value = "  -0.05%  "

# Row 19
$ws.Range("D19").Value = "'432.52"
$ws.Range("E19").Value = "  +0.21%  "

# Row 20
$ws.Range("D20").Value = "'13.49"
$ws.Range("E20").Value = "  +0.30%  "

# Row 21
$ws.Range("D21").Value = "'0.678"
$ws.Range("E21").Value = "  -0.81%  "

# Row 22
$ws.Range("D22").Value = "'7.11"
$ws.Range("E22").Value = "  -0.01%  "

# Row 23
$ws.Range("D23").Value = "'81.93"
$ws.Range("E23").Value = "  +0.81%  "

# Row 24
$ws.Range("E24").Value = "  -0.85%  "

# Row 25
$ws.Range("D25").Value = "'2.20"
$ws.Range("E25").Value = "  -1.36%  "

# Row 26
$ws.Range("E26").Value = "  -1.91%  "

# Row 28
$ws.Range("E28").Value = "  -3.47%  "

# Row 29
$ws.Range("E29").Value = "  -0.11%  "

# Row 30
$ws.Range("D30").Value = "'6.95"
$ws.Range("E30").Value = "  -2.12%  "

# Row 31
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.110"
$ws.Range("E31").Value = "  +1.70%  "

# Row 32
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").Value = "'26.70"
$ws.Range("E32").Value = "  +0.73%  "

# Row 33
$ws.Range("E33").Value = "  +0.05%  "

# Row 34
$ws.Range("E34").Value = "  +2.95%  "

# Row 35
$ws.Range("E35").Value = "  +0.12%  "

# Row 36
$ws.Range("D36").Value = "'5.64"
$ws.Range("E36").Value = "  +0.41%  "

# Row 37
$ws.Range("D37").Value = "'3.00"
$ws.Range("E37").Value = "  -2.21%  "

# Row 38
$ws.Range("E38").Value = "  +0.64%  "

# Row 39
$ws.Range("D39").Value = "'0.123"
$ws.Range("E39").Value = "  -1.02%  "

# Row 40
$ws.Range("D40").Value = "'8.59"
$ws.Range("E40").Value = "  -0.03%  "

# Row 41
$ws.Range("D41").Value = "'42.37"
$ws.Range("E41").Value = "  +7.59%  "

# Row 42
$ws.Range("D42").Value = "'0.283"
$ws.Range("E42").Value = "  -1.39%  "

# Row 43
$ws.Range("D43").Value = "'0.0347"
$ws.Range("E43").Value = "  +0.04%  "

# Row 44
$ws.Range("D44").Value = "2.701.53"
$ws.Range("E44").Value = "  -0.46%  "

# Row 45
$ws.Range("D45").Value = "'134.39"
$ws.Range("E45").Value = "  +1.71%  "

# Row 46
$ws.Range("D46").Value = "'364.80"
$ws.Range("E46").Value = "  -2.87%  "

# Row 48
$ws.Range("D48").Value = "'23.79"
$ws.Range("E48").Value = "  -1.49%  "

# Row 49
$ws.Range("E49").Value = "  -1.34%  "

# Row 50
$ws.Range("E50").Value = "  -1.39%  "

# Row 51
$ws.Range("D51").Value = "'0.124"
$ws.Range("E51").Value = "  -1.35%  "
